# Added scroll wheel for channels in configure time mode.
# - Typography sheet: existing Default/Large/Small fonts get a "0-9"
#   wildcard range, and two new typography entries (scrollWheel,
#   scrollWheelCenter) are added using Roboto-Bold.ttf.
# - Translation sheet: new text ids wiring up the new typography names.

$wb = $excel.ActiveWorkbook

$typography = $wb.Worksheets.Item("Typography")
$translation = $wb.Worksheets.Item("Translation")

# --- Typography sheet -------------------------------------------------
# Existing rows 4-6 (Default/Large/Small) gain a "0-9" wildcard range.
$typography.Range("I4").Value = "0-9"
$typography.Range("I5").Value = "0-9"
$typography.Range("I6").Value = "0-9"

# New row 7: scrollWheel typography entry
$typography.Range("B7").Value = "scrollWheel"
$typography.Range("B7").Style = "Normal"
$typography.Range("C7").Value = "Roboto-Bold.ttf"
$typography.Range("C7").Style = "Normal"
$typography.Range("D7").Value = 25
$typography.Range("D7").Style = "Normal"
$typography.Range("E7").Value = 4
$typography.Range("E7").Style = "Normal"
$typography.Range("F7").Value = "?"
$typography.Range("G7").Value = "'"
$typography.Range("G7").Style = "Normal"
$typography.Range("H7").Value = "'"
$typography.Range("H7").Style = "Normal"
$typography.Range("I7").Value = "0-9"
$typography.Range("J7").Value = "'"
$typography.Range("J7").Style = "Normal"

# New row 8: scrollWheelCenter typography entry
$typography.Range("B8").Value = "scrollWheelCenter"
$typography.Range("B8").Style = "Normal"
$typography.Range("C8").Value = "Roboto-Bold.ttf"
$typography.Range("C8").Style = "Normal"
$typography.Range("D8").Value = 35
$typography.Range("D8").Style = "Normal"
$typography.Range("E8").Value = 4
$typography.Range("E8").Style = "Normal"
$typography.Range("F8").Value = "?"
$typography.Range("G8").Value = "'"
$typography.Range("G8").Style = "Normal"
$typography.Range("H8").Value = "'"
$typography.Range("H8").Style = "Normal"
$typography.Range("I8").Value = "0-9"
$typography.Range("J8").Value = "'"
$typography.Range("J8").Style = "Normal"

# --- Translation sheet -------------------------------------------------
# New rows 17-21 wiring the new typography names to text ids.
$translation.Range("B17").Value = "SingleUseId16"
$translation.Range("C17").Value = "scrollWheel"
$translation.Range("D17").Value = "Center"
$translation.Range("E17").Value = "LTR"
$translation.Range("F17").Value = "<value>"

$translation.Range("B18").Value = "SingleUseId17"
$translation.Range("C18").Value = "scrollWheel"
$translation.Range("D18").Value = "Left"
$translation.Range("E18").Value = "LTR"
$translation.Range("F18").Value = "'1"
$translation.Range("F18").Style = "Normal"

$translation.Range("B19").Value = "SingleUseId14"
$translation.Range("C19").Value = "scrollWheelCenter"
$translation.Range("D19").Value = "Center"
$translation.Range("E19").Value = "LTR"
$translation.Range("F19").Value = "<value>"

$translation.Range("B20").Value = "SingleUseId15"
$translation.Range("C20").Value = "scrollWheelCenter"
$translation.Range("D20").Value = "Left"
$translation.Range("E20").Value = "LTR"
$translation.Range("F20").Value = "'1"
$translation.Range("F20").Style = "Normal"

$translation.Range("B21").Value = "SingleUseId18"
$translation.Range("C21").Value = "Default"
$translation.Range("D21").Value = "Center"
$translation.Range("E21").Value = "LTR"
$translation.Range("F21").Value = "<value>"
